$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 13 play-by-play yardage logs to the running totals ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = "6 8 0 -2 15 2 2 -2 5 3 6 0 2 5 12 2 3 3 0 4 2 2 0 16 2 5 -2 13 -1 4 5 1 11 3 10 1 4 5 1 3 2 4 14 0 5 1 3 4 14 4 3 0 6 1 7 15 22 3 1 2 2 44 5 3 30 5 6 11 3 6 4 18 4 8 11 1 4 5 24 3 4 1 8 1 2 5 1 1 5 1 50 1 2 7 2 3 1 3 4 0 13 8 6 9 11 -2 0 14 5 42 3 4 1 1 3 -4 4 0 6 5 3 34 0 3 10 3 14 3 20 4 3 3 0 0 1 -2 9 4 -2 2 3 7 2 7 5 -2 0 10 16 5 6 2 1 -2 24 -1 4 6 1 37 0 5 -4 3 1 10 5 21 4 4 -1 1 4 3 4 3 28 1 6 5 25 6 12 3 14 3 5 1 4 2 0 4 8 9 4 8 2 4 1 4 2 0 3 -2 3 7 4 9 5 20 15 3 2 2 3 -1 3 -2 7 -1 1 8 -2 5 2 8 4 4 3 3 4 5 4 1 0 1 5 5 1 1 3 8 5 8 1 9 4 5 -1 3 5 -2 -1 4 4 3 1 5 7 5 3 3 9 5 1 3 5 6 1 2 1 8 5 -2 7 11 0 8 3 2 5 0 2 10 7 1 2 6 9 -3 1 11 4 7 2 2 1 3 2 3 2 3 2 5 0 1 8 8 14 2 1 -4 1 5 6 11 3 1 -1 -1 39 4 4 1 19 3 0 6 0 5 7 -2 -2 6 0 1 9 3 10 7 5 18 0 1 36 4 7 2 0 2 30 -1 4 2 2 7 7 5 5 14 21 9 -2 4 6 5 0 7 0 5 24 3 5 5 7 4 13 6 1 5 11 4 10 13 1 20 11 2 4 17 1 44 19 3 1 2 4 1 4 0 5 8 6 2 -3 3 4 11 7 1 8 2 9 3 5 2 3 -3 5 0 7 6 3 0 4 3 18 5 4 0 7 2 2 6 19 1 17 5 4 3 5 6 2 7 8 5 1 18 3 6 2 6 17 2 9 5 7 3 2 7 7 20 32 0 7 1 25 5 6 3 6 0 3 0 7 0 0 4 3 1 13 7 2 22 2 12 3 5 6 3 6 5 5 5 8 22 9 0 3 2 5 3 2 3 18 6 27 20 4 4 7 4 4 11 4 8 72 16 6 1 16 2 2 5 2 0 1 2 0 1 3 2 4 0 3 3 5 0 4 8 1 15 2 48 8 10 -1 9 23 2 4 4 -2 1 3 12 2 14 4 33 5 4 5 12 10 9 -1 5 7 1 2 5 0 3 3 0 1 2 13 -4 9 5 2 15 1 1 6 -4 0 19 3 9 2 11 -2 3 18 2 2 3 -3 0 4 1 16 3 2 4 11 6 -1 15 11 4 2 2 5 5 8 1 0 28 2 4 2 3 2 6 5 12 31 9 7 1 9 9 8 0 0 21 3 4 5 20 14 -4 3 4 -2 20 6 2 8 5 9 2 1 11 5 9 -1 3 1 1 2 2 -2 2 5 0 31 6 2 0 4 8 8 1 9 0 5 7 8 5 12 3 -1 7 5 2 1 1 3 4 5 7 11 2 1 -1 6 -2 3 4 3 5 4 0 2 3 5 3 5 5 4 5 0 1 3 3 2 -2 13 4 3 -1 9 2 4 0 12 5 -1 3 6 4 -3 6 6 6 2 8 6 16 6 1 11 14 3 3 22 2 3 5 5 7 3 6 1 5 4 1 7 2 2 0 -3 5 7 3 0 9 -1 9 3 2 12 2 10 -5 4 -2 -1 12 3 10 3 7 1 16 2 1 7 7 2 10 7 2 5 6 5 8 5 6 4 4 2 6 3 4 9 3 5 6 7 11 6 12 0 1 7 1 4 9 4 6 1 3 8 7 8 5 10 1 6 5 4 2 3 19 5 12 1 4 19 4 12 2 8 4 2 -2 0 0 5 2 11 0 -4 1 3 1 11 5 7 5 10 1 -1 -3 4 1 4 2 5 4 2 19 0 3 8 2 3 4 5 2 5 4 6 5 0 2 1 6 -1 2 2 3 1 2 4 1 1 6 1 0 3 4 0 13 9 5 4 2 0 5 1 0 4 2 4 1 4 3 2 8 13 1 3 11 8 3 1 4 6 5 4 1 12 5 18 -1 3 5 3 8 -2 4 5 3 8 6 8 0 -1 -2 13 3 2 2 3 1"
$wsYDS.Range("C2").Value = "11 0 4 2 29 22 1 3 8 3 0 0 0 4 1 3 2 10 10 9 4 3 4 3 1 4 0 1 1 8 2 5 2 5 6 9 1 -2 4 2 2 22 1 3 4 5 1 1 -1 9 4 16 1 7 3 11 2 0 -1 0 12 6 5 3 3 1 2 4 2 2 5 1 10 2 2 -1 3 7 -1 3 9 1 -1 2 9 5 3 7 3 2 0 1 5 5 1 6 3 1 -3 2 2 2 2 6 3 1 2 -1 4 3 7 2 0 7 2 0 2 2 -1 0 2 4 2 5 1 20 3 1 26 5 74 5 -1 40 5 1 -1 1 9 4 4 0 18 6 -2 3 3 1 -3 1 1 3 0 11 3 4 0 4 16 2 1 9 12 3 1 2 4 3 11 5 4 2 1 14 7 4 7 3 0 6 7 4 3 4 4 12 5 1 13 4 16 25 4 12 3 4 1 -2 -1 0 9 1 1 0 -1 8 8 1 2 -2 4 -5 5 2 4 8 8 3 0 21 3 3 3 3 0 2 5 1 0 2 -2 24 1 11 3 1 12 4 3 8 4 4 29 5 2 1 3 3 -2 7 3 1 2 3 10 3 3 2 13 7 2 1 1 7 2 9 2 13 4 1 2 10 3 -1 0 10 14 3 2 2 3 2 3 3 0 -1 0 3 12 3 3 7 2 2 7 4 2 3 8 14 1 11 5 4 5 -1 2 10 2 4 5 7 3 6 13 1 2 2 10 3 0 6 -1 2 1 0 -1 2 0 2 0 1 2 2 11 2 -1 17 3 5 6 1 1 9 3 3 3 15 3 2 1 3 4 13 1 11 2 -2 55 2 -2 4 6 15 9 2 4 17 4 2 6 1 1 -3 5 4 1 1 0 4 2 4 1 3 2 0 8 3 0 4 3 9 -1 12 6 3 -2 -1 0 3 1 31 -1 2 8 0 13 0 0 0 2 2 2 4 1 15 3 -2 6 0 1 9 1 7 -1 1 3 9 4 2 7 6 15 1 3 -5 4 -2 -2 11 1 4 14 2 6 4 14 4 0 9 -2 2 1 2 3 12 4 0 4 0 1 1 -1 0 1 -4 7 8 1 4 14 3 2 31 3 3 3 3 -3 6 10 -1 12 3 11 4 4 3 13 2 2 -1 2 2 19 13 0 4 4 8 22 -1 7 2 1 -1 3 -4 1 -1 4 2 5 2 10 -2 2 1 3 -2 2 9 0 -1 -1 9 0 0 5 1 0 0 10 7 21 3 46 3 -3 -5 2 2 1 11 15 -3 4 66 1 3 3 0 0 1 2 -1 1 0 -1 0 -8 1 24 1 3 -1 3 1 -1 5 4 2 4 2 1 -1 4 1 14 0 11 0 1 7 2 1 8 9 2 11 5 -6 3 10 7 5 6 4 5 5 -1 6 11 2 0 -3 4 1 2 7 2 2 -1 3 4 -1 5 1 9 3 -3 5 2 7 0 5 -1 5 4 2 3 6 2 2 0 1 8 3 2 0 8 13 3 5 8 -2"
$wsYDS.Range("B3").Value = "12 47 9 25 21 4 16 20 13 16 11 9 4 21 19 12 -1 2 7 10 9 6 -3 4 29 14 4 18 11 1 19 25 11 7 13 5 13 22 5 7 0 0 7 4 6 13 7 7 8 5 19 6 6 4 10 1 33 31 3 25 13 9 18 22 4 4 10 10 7 16 9 9 14 5 14 2 7 2 9 12 9 12 3 7 5 8 30 11 17 7 9 32 9 11 5 18 2 -1 20 10 17 10 9 14 34 6 3 11 6 10 18 39 24 3 8 4 6 4 5 4 7 13 9 -1 14 14 9 21 20 12 4 11 9 5 17 13 21 19 6 11 31 2 3 8 9 9 12 6 7 4 14 8 6 18 10 11 3 1 5 14 30 5 4 6 6 31 14 5 3 6 8 15 16 9 9 5 -1 5 6 7 3 8 7 70 1 3 4 12 8 4 38 4 2 20 9 2 4 10 39 6 15 7 13 14 14 6 4 11 10 18 28 4 27 -5 11 2 11 6 6 34 -1 44 12 3 -1 8 13 13 9 6 11 6 11 20 0 7 25 6 10 4 8 27 7 5 10 43 5 9 18 5 2 9 8 7 17 5 28 9 19 -6 11 4 6 9 3 17 20 9 17 4 1 11 31 11 21 8 30 5 5 7 14 1 8 9 5 2 4 9 11 29 2 8 0 21 6 6 29 29 10 13 4 10 6 4 13 6 49 9 10 11 17 -4 7 20 8 9 7 11 4 4 12 8 27 20 42 19 9 10 6 5 16 11 6 21 37 19 29 3 9 41 7 24 5 18 20 6 12 24 12 9 49 32 7 13 14 4 8 13 7 8 17 4 7 14 32 4 5 5 17 15 28 7 28 2 6 11 12 8 8 7 17 3 19 7 43 27 18 7 6 15 5 17 10 9 20 8 4 14 4 11 7 7 5 7 4 4 8 16 8 9 21 1 10 6 4 11 9 8 20 13 6 1 24 6 23 19 20 18 7 1 35 39 25 3 13 13 -1 10 12 7 11 13 13 15 7 -1 9 4 5 12 22 3 3 8 8 18 11 5 22 6 8 7 7 20 20 7 13 7 4 11 15 7 5 6 13 3 11 9 9 11 0 8 7 8 4 8 20 5 8 30 6 9 15 11 1 16 5 9 8 -2 11 11 2 3 12 9 15 2 9 8 -2 3 9 4 12 7 21 6 29 8 2 7 5 6 9 8 13 8 8 4 7 39 13 11 5 -1 0 10 3 14 16 18 0 9 9 8 2 10 29 3 10 15 15 3 17 6 5 22 18 13 5 6"
$wsYDS.Range("C3").Value = "5 6 16 6 1 8 13 9 1 16 28 8 5 21 3 -2 4 0 9 12 4 7 13 6 6 13 38 10 6 5 12 9 13 15 9 6 6 15 4 24 11 7 11 23 3 11 33 15 11 13 15 2 5 7 29 5 18 -4 20 10 49 9 5 7 24 7 6 9 7 6 18 13 11 2 17 5 9 17 -2 3 11 8 18 11 1 13 -1 40 4 3 12 0 7 4 4 19 4 13 8 8 11 11 9 16 4 39 14 26 12 -2 14 6 7 10 4 16 8 11 3 5 12 9 9 7 -6 9 8 3 19 13 8 6 4 50 6 11 3 5 9 10 10 10 1 16 10 0 7 24 18 14 9 8 13 7 14 0 9 6 4 5 2 11 8 6 2 3 13 20 13 17 5 12 10 17 5 9 1 5 -2 3 18 6 6 5 6 4 21 3 4 19 20 7 4 4 24 26 8 6 12 1 4 14 7 8 1 -1 7 6 2 12 50 5 12 16 9 5 9 22 25 14 14 11 11 13 1 4 5 3 17 27 6 11 5 4 7 22 8 13 -1 5 11 4 1 5 6 10 8 7 8 1 7 7 5 8 5 3 3 3 16 4 9 -4 0 2 14 13 11 13 12 13 6 6 15 9 4 6 13 13 9 11 6 4 19 19 12 8 12 -4 20 10 12 4 11 9 4 4 11 37 7 18 1 9 14 9 26 6 5 12 8 10 30 7 16 22 8 8 -2 12 10 2 3 12 11 21 8 6 17 22 2 28 10 8 9 4 13 10 7 14 13 9 11 11 4 15 16 19 20 8 -5 7 0 9 16 13 11 3 12 14 15 1 8 21 8 7 3 28 4 14 5 9 10 35 5 10 6 7 18 12 3 -2 10 6 5 13 4 16 1 11 21 2 11 10 2 4 19 6 14 20 9 3 4 7 14 10 1 4 9 13 6 5 6 12 9 24 8 15 13 7 4 12 16 2 3 4 21 4 5 16 9 7 9 37 10 10 20 18 11 27 32 31 6 7 19 10 33 17 11 14 10 4 14 20 14 5 11 40 5 46 3 8 3 23 13 7 11 3 4 5 8 1 15 10 5 12 11 14 19 9 2 -3 11 9 6 24 19 22 11 19 3 3 12 14 3 11 11 8 3 5 6 13 -3 9 9 32 9 76 11 15 6 8 15 11 7 -2 5 9 16 -1 31 42 21 24 18 4 4 29 23 2 15 11 11 8 3 18 7 5 8 26 1 16 1 5 10 21 15 -1 5 4 5 10 6 11 23 6 5 9 55 7 4 6 13 26 7 17 21 32 7 82 25 2 11 11 27 6 50 5 3 7 7 20 4 18 5 9 5 6 27 1 2 16 9 11 7 21 2 6 19 2 52 11 9 1 6 14 35 23 6 64 6 29 22 22 23 60 10 2 3 12 9 4 23 11 9 36 10 0 -1 41 11 7 4 16 6 38 8 20 16 8 3 18 7 12 4 17 5 14 7 7 4 9 40 29 25 9 2 5 9 5 5"

# --- OFF sheet: updated cumulative totals ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 434
$wsOFF.Range("D2").Value = 21
$wsOFF.Range("E2").Value = 23
$wsOFF.Range("F2").Value = 128
$wsOFF.Range("G2").Value = 133
$wsOFF.Range("H2").Value = 11
$wsOFF.Range("J2").Value = 104
$wsOFF.Range("N2").Value = 55
$wsOFF.Range("B3").Value = 16
$wsOFF.Range("C3").Value = 284
$wsOFF.Range("E3").Value = 57
$wsOFF.Range("F3").Value = 190
$wsOFF.Range("G3").Value = 37
$wsOFF.Range("H3").Value = 52
$wsOFF.Range("I3").Value = 105
$wsOFF.Range("J3").Value = 56
$wsOFF.Range("L3").Value = 459
$wsOFF.Range("M3").Value = 293
$wsOFF.Range("Q3").Value = 1033

# --- DEF sheet: updated cumulative totals ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 341
$wsDEF.Range("F2").Value = 87
$wsDEF.Range("G2").Value = 80
$wsDEF.Range("J2").Value = 36
$wsDEF.Range("N2").Value = 35
$wsDEF.Range("C3").Value = 324
$wsDEF.Range("F3").Value = 206
$wsDEF.Range("G3").Value = 61
$wsDEF.Range("H3").Value = 53
$wsDEF.Range("I3").Value = 111
$wsDEF.Range("J3").Value = 94
$wsDEF.Range("L3").Value = 585
$wsDEF.Range("M3").Value = 358
$wsDEF.Range("Q3").Value = 965

# --- ST sheet: updated cumulative totals + appended per-kick/punt logs ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 165
$wsST.Range("D2").Value = 113
$wsST.Range("F2").Value = 377
$wsST.Range("G2").Value = 373
$wsST.Range("H2").Value = 8
$wsST.Range("J2").Value = 177
$wsST.Range("K2").Value = 174
$wsST.Range("B3").Value = 108
$wsST.Range("B4").Value = "65 64 62 60 61 61 65 67 55 63 68 70 63 63 59 65 62 58 61 59 68 65 62 55 62 65 61 59 65 61 66 65 57 63 65 57 65 55 48 67 58 63 64 65 61 63 60 67 59 67 63 61 39 65 61 64 58"
$wsST.Range("B5").Value = "17 19 10 20 22 21 27 22 12 24 30 37 14 14 19 20 26 6 15 66 29 22 16 19 24 23 19 18 18 16 25 21 26 13 20 8 33 3 9 25 17 19 26 20 16 18 16 36 10 23 98 17 0 15 14 19 13"
$wsST.Range("B6").Value = "38 26 21 20 16 26 93 30 37 20 26 24 37 17 26 18 28 26 24 21 5 16 23 20 23 23 21 23 27 23 47 21 24 24 21 25 25 26 35 0"
$wsST.Range("D3").Value = "49 37 39 53 41 48 49 50 49 41 51 44 42 42 49 57 59 45 46 26 43 49 36 44 39 50 43 40 44 43 48 44 58 49 40 40 46 42 58 46 35 52 46 42 34 46 50 35 41 28 40 55 45 53 43 43 48 23 37 48 55 47 42 49 41 49 46 54 48 35 34 52 55 46 49 57 42 49 53 47 49 53 33 48 49 40 31 50 46 39 35 51 45 37 53 43 39 24 41 39 37 42 44 35 57 45 53 33 50 55 45 46 57"
$wsST.Range("D4").Value = "1 0 0 10 0 0 7 0 -3 0 0 4 0 0 0 11 5 0 0 0 0 4 5 0 0 11 9 0 0 -2 3 0 0 0 0 0 0 0 0 0 0 6 0 0 0 0 5 0 0 0 6 28 -2 6 0 0 2 0 0 8 0 3 0 0 7 0 13 11 0 0 0 0 5 8 0 9 0 6 0 13 11 0 0 13 7 0 0 0 0 0 0 0 6 0 0 0 0 0 4 0 0 0 0 0 4 0 0 0 -3 14 5 6 8"
$wsST.Range("D5").Value = "0 17 9 0 0 8 16 0 10 5 13 0 0 0 0 10 1 -1 0 10 0 3 5 10 11 9 0 0 0 0 15 0 0 0 9 4 0 0 0 0 0 7 0 11 0 0 0 0 0 0 20 32 17 19 0 14 0 0 0 6 0 0 0 0 0 0 0 0 0 0 0 0 22 15 30 0 0 16 0 0 0 0 13 29 0 12 0 0 10 0 42 0 0 0 0 0 7 0 14 13 17 0 0 16 21 0 0 0 0 14 0 0 0 0 0 9 0 0 0 7 12 0 0 0 13 22 0 0 1 0 0 0 12 0 0"

# --- TURNS sheet: updated cumulative totals ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 12
$wsTURNS.Range("C3").Value = 10
$wsTURNS.Range("D3").Value = 20
$wsTURNS.Range("E3").Value = 17

# --- PEN sheet: updated cumulative totals ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 35
$wsPEN.Range("D2").Value = 12
$wsPEN.Range("B3").Value = 41
$wsPEN.Range("D4").Value = 22
